$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.817.63"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.313.85"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").Value = "2.676.01"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "2.341.60"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "42.764.27"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +5.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +14.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0698"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.60%  "
$ws.Range("D43").Value = "1.927.50"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "2.544.45"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
